# Update the "dSF" column (column F) values for a set of rows.
# These values are being re-pulled / recalculated from source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 7
    5  = -2
    7  = -3
    8  = -11
    13 = -3
    23 = -17
    31 = -8
    39 = -8
    41 = -6
    42 = 7
    43 = -6
    46 = -4
    50 = -5
    51 = -6
    55 = -3
    56 = 2
    58 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
